$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.8125
$ws.Range("C2").Value = 0.7959183673469388
$ws.Range("D2").Value = 0.8041237113402061
$ws.Range("E2").Value = 49
$ws.Range("B3").Value = 0.7777777777777778
$ws.Range("C3").Value = 0.7954545454545454
$ws.Range("D3").Value = 0.7865168539325843
$ws.Range("E3").Value = 44
$ws.Range("B4").Value = 0.7956989247311828
$ws.Range("C4").Value = 0.7956989247311828
$ws.Range("D4").Value = 0.7956989247311828
$ws.Range("E4").Value = 0.7956989247311828
$ws.Range("B5").Value = 0.7951388888888888
$ws.Range("C5").Value = 0.7956864564007421
$ws.Range("D5").Value = 0.7953202826363952
$ws.Range("B6").Value = 0.7960722819593788
$ws.Range("C6").Value = 0.7956989247311828
$ws.Range("D6").Value = 0.7957935852548795
$ws.Range("B7").Value = 0.7916666666666666
$ws.Range("C7").Value = 0.7755102040816326
$ws.Range("D7").Value = 0.7835051546391752
$ws.Range("E7").Value = 49
$ws.Range("B8").Value = 0.7555555555555555
$ws.Range("C8").Value = 0.7727272727272727
$ws.Range("D8").Value = 0.7640449438202247
$ws.Range("E8").Value = 44
$ws.Range("B9").Value = 0.7741935483870968
$ws.Range("C9").Value = 0.7741935483870968
$ws.Range("D9").Value = 0.7741935483870968
$ws.Range("E9").Value = 0.7741935483870968
$ws.Range("B10").Value = 0.773611111111111
$ws.Range("C10").Value = 0.7741187384044527
$ws.Range("D10").Value = 0.7737750492297
$ws.Range("B11").Value = 0.7745818399044204
$ws.Range("C11").Value = 0.7741935483870968
$ws.Range("D11").Value = 0.774298173176446
$ws.Range("B12").Value = 0.78
$ws.Range("C12").Value = 0.7959183673469388
$ws.Range("D12").Value = 0.7878787878787878
$ws.Range("E12").Value = 49
$ws.Range("B13").Value = 0.7674418604651163
$ws.Range("C13").Value = 0.75
$ws.Range("D13").Value = 0.7586206896551724
$ws.Range("E13").Value = 44
$ws.Range("B14").Value = 0.7741935483870968
$ws.Range("C14").Value = 0.7741935483870968
$ws.Range("D14").Value = 0.7741935483870968
$ws.Range("E14").Value = 0.7741935483870968
$ws.Range("B15").Value = 0.7737209302325582
$ws.Range("C15").Value = 0.7729591836734694
$ws.Range("D15").Value = 0.7732497387669801
$ws.Range("B16").Value = 0.7740585146286572
$ws.Range("C16").Value = 0.7741935483870968
$ws.Range("D16").Value = 0.7740362467837439
$ws.Range("B17").Value = 0.8627450980392157
$ws.Range("C17").Value = 0.8979591836734694
$ws.Range("D17").Value = 0.8799999999999999
$ws.Range("E17").Value = 49
$ws.Range("B18").Value = 0.8809523809523809
$ws.Range("C18").Value = 0.8409090909090909
$ws.Range("D18").Value = 0.8604651162790699
$ws.Range("E18").Value = 44
$ws.Range("B19").Value = 0.8709677419354839
$ws.Range("C19").Value = 0.8709677419354839
$ws.Range("D19").Value = 0.8709677419354839
$ws.Range("E19").Value = 0.8709677419354839
$ws.Range("B20").Value = 0.8718487394957983
$ws.Range("C20").Value = 0.8694341372912802
$ws.Range("D20").Value = 0.8702325581395349
$ws.Range("B21").Value = 0.8713592964067348
$ws.Range("C21").Value = 0.8709677419354839
$ws.Range("D21").Value = 0.8707576894223557
$ws.Range("B22").Value = 0.8222222222222222
$ws.Range("C22").Value = 0.7551020408163265
$ws.Range("D22").Value = 0.7872340425531914
$ws.Range("E22").Value = 49
$ws.Range("B23").Value = 0.75
$ws.Range("C23").Value = 0.8181818181818182
$ws.Range("D23").Value = 0.7826086956521738
$ws.Range("E23").Value = 44
$ws.Range("B25").Value = 0.7861111111111111
$ws.Range("C25").Value = 0.7866419294990723
$ws.Range("D25").Value = 0.7849213691026826
$ws.Range("B26").Value = 0.78805256869773
$ws.Range("D26").Value = 0.7850457063849681
